# Rename the inline picture shapes in the document's headers/footers:
#   - the two Pearson logo pictures (currently "image1.png") become "image2.png"
#   - the BTEC logo picture (currently "image2.jpg") becomes "image1.jpg"
#
# InlineShape has no writable Name property (matches real Word's object
# model), so each picture is temporarily converted to a floating Shape -
# which does expose Name - renamed, then converted back to an inline
# shape so the picture's layout/wrapping is left exactly as it was.

$d = $word.ActiveDocument

function Rename-InlinePicture($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le 3; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $pic = $shapes.Item($i)
                if ($pic.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlinePicture $pic "image1.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le 3; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $pic = $shapes.Item($i)
                if ($pic.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    Rename-InlinePicture $pic "image2.png"
                }
            }
        }
    }
}
